# edit.ps1 - apply the four run-level changes described by the diff.
#
# 1) "Création du design ou choix du template " -> split "template" into
#    its own run (originally wrapped in <w:proofErr> spell-check markers
#    by Word's live proofer; those markers are not reachable through the
#    Word COM surface, so only the run split itself is reproduced).
# 2) "Création des Models" -> split "Models" into its own run (same
#    proofErr caveat as above).
# 3) "... avec envoi de mail" -> append a space run + a bold red "…" run.
# 4) "... avec ajout d'un nouveau profil" -> append a space run + a bold
#    red "-bof" run.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: split "template" out of its sentence.
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Création du design ou choix du template", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $sentenceStart = $rng.Start
    $wordStart = $sentenceStart + 31
    $wordEnd = $wordStart + 8
    $tmplRange = $d.Range($wordStart, $wordEnd)
    # Force Word to materialize this span as its own run (toggling a real
    # boolean property and then reverting it is the only reliable way to
    # split an otherwise-identical run in this runtime).
    $tmplRange.Bold = $true
    $tmplRange.Bold = $false
}

# ---------------------------------------------------------------------
# Change 2: split "Models" out of its sentence.
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Création des Models", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $sentenceStart = $rng.Start
    $wordStart = $sentenceStart + 13
    $wordEnd = $wordStart + 6
    $modelsRange = $d.Range($wordStart, $wordEnd)
    $modelsRange.Bold = $true
    $modelsRange.Bold = $false
}

# ---------------------------------------------------------------------
# Change 3: after "... avec envoi de mail" add " " + bold red "…"
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("avec envoi de mail", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Collapse(0)
    $rng.InsertAfter(" ")
    $rng.Collapse(0)
    $markStart = $rng.End
    $rng.InsertAfter("…")
    $markRange = $d.Range($markStart, $rng.End)
    $markRange.Font.Bold = $true
    $markRange.Font.Color = 255
}

# ---------------------------------------------------------------------
# Change 4: after "... avec ajout d'un nouveau profil" add " " + bold
# red "-bof"
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("avec ajout d" + [char]0x2019 + "un nouveau profil", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Collapse(0)
    $rng.InsertAfter(" ")
    $rng.Collapse(0)
    $markStart = $rng.End
    $rng.InsertAfter("-bof")
    $markRange = $d.Range($markStart, $rng.End)
    $markRange.Font.Bold = $true
    $markRange.Font.Color = 255
}

Write-Host "done"
